# Atualização automática de preços de eletricidade
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45980
$ws.Range("B2").Value = 119.36
$ws.Range("C2").Value = 104.93
$ws.Range("D2").Value = 97.01000000000001
$ws.Range("E2").Value = 96.91
$ws.Range("F2").Value = 95.59
$ws.Range("G2").Value = 97.95999999999999
$ws.Range("H2").Value = 113.23
$ws.Range("I2").Value = 114.14
$ws.Range("J2").Value = 116.62
$ws.Range("K2").Value = 99.45
$ws.Range("L2").Value = 92.11
$ws.Range("M2").Value = 74.91
$ws.Range("N2").Value = 70.95999999999999
$ws.Range("O2").Value = 66.98999999999999
$ws.Range("P2").Value = 69.06
$ws.Range("Q2").Value = 82.22
$ws.Range("R2").Value = 89.59999999999999
$ws.Range("S2").Value = 114.24
$ws.Range("T2").Value = 120.02
$ws.Range("U2").Value = 124.06
$ws.Range("V2").Value = 129.25
$ws.Range("W2").Value = 114.51
$ws.Range("X2").Value = 104.16
$ws.Range("Y2").Value = 96.68000000000001
$ws.Range("Z2").Value = 100.17
$ws.Range("AA2").Value = "16h-20h"
$ws.Range("AB2").Value = 111.98
$ws.Range("AC2").Value = "18h-20h"
$ws.Range("AD2").Value = 122.04
$ws.Range("AE2").Value = "20h-22h"
$ws.Range("AF2").Value = 121.88
$ws.Range("AG2").Value = "2h-23h"
